# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were recomputed from the newly regenerated
# save_data (std/mean) and the resulting s_vals are written back into column G
# for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 6
    5  = 2
    6  = 2
    7  = 0
    8  = 3
    9  = 2
    10 = 1
    11 = 1
    12 = 3
    13 = 7
    14 = 1
    15 = 3
    16 = 3
    17 = 4
    18 = 2
    19 = 6
    20 = 5
    21 = 1
    22 = 3
    23 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
